# Implemented new methods on ExcelPackage class, and renamed some methods
#
# - TEST4: drop the "R28" selection remnant, select the real table range
#   (E9:G13) instead, and size columns E/G to fit their (header/date)
#   content.
# - TEST5: clear the stray "H14" selection, just mark the used range
#   (A1:C4) as selected.
# - Add a new sheet "TEST6" (a copy of TEST5's layout/data) with the last
#   quantity cell left blank, and make it the active sheet/selection.

$wb = $excel.ActiveWorkbook

# ---- TEST4 -----------------------------------------------------------
$ws4 = $wb.Worksheets.Item("TEST4")
$ws4.Range("E9:G13").Select()
$ws4.Columns.Item(5).ColumnWidth = 21.166666666666668   # -> width 22
$ws4.Columns.Item(7).ColumnWidth = 14.307291666666666   # -> width ~15.14

# ---- TEST5 -------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("TEST5")
$ws5.Range("A1:C4").Select()

# ---- TEST6 (new sheet, appended after TEST5) ---------------------------
$ws6 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws6.Name = "TEST6"

$ws5.Range("A1:C4").Copy($ws6.Range("A1"))
$ws6.Range("B4").ClearContents()
$ws6.Range("C4").ClearContents()

$ws6.Range("H17").Select()
